$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '245.57'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '23.94'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.349'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05817'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.365'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.465'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8096'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9217'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.01069'
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1401'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07402'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03195'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03063'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09363'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.847'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.001548'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.04708'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006037'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.001244'
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.004695'
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.00008794'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B23').Value = 'LEO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.597'
$ws.Range('E23').Value = '22LEOLEO'
$ws.Range('B24').Value = 'BTSEToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.150'
$ws.Range('E24').Value = '23BTSETokenBTSE'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.3184'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0002348'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006376'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1067'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005247'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000749'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6850'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.001828'
